$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2: "Clase gastronomica" (virtual, artistic class) ---
$ws.Range("A2").Value = "Clase gastronomica"
$ws.Range("B2").Value = "artístico"
$ws.Range("C2").Value = "virtual"
$ws.Range("D2").Value = "Gastronomia virtual"
$ws.Range("E2").Value = "no aplica"
$ws.Range("F2").Value = "virtual"
$ws.Range("G2").Value = 3118522584
$ws.Range("H2").Value = "No aplican requisitos"
$ws.Range("I2").Value = 120

# --- Update row 3: "Clase astronomica" (virtual, artistic class) ---
$ws.Range("A3").Value = "Clase astronomica"
$ws.Range("B3").Value = "artístico"
$ws.Range("C3").Value = "virtual"
$ws.Range("D3").Value = "Gastronomia virtual"
$ws.Range("E3").Value = "no aplica"
$ws.Range("F3").Value = "virtual"
$ws.Range("G3").Value = 3118522584
$ws.Range("H3").Value = "No aplican requisitos"
$ws.Range("I3").Value = 80

# --- Remove the remaining sample rows (old rows 4-7), replacing them with
#     fresh blank rows so the leftover custom row height (30pt) used by the
#     sample data rows is dropped and the sheet returns to its template
#     layout (blank rows 4-6, then the plain G/H placeholder pattern
#     resuming at row 7, same as the rows below it). ---
$ws.Range("4:7").Delete()
$ws.Range("4:6").Insert()

# --- Update the saved selection to reflect where editing left off ---
[void]$ws.Range("I7").Select()
